$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the B column (path) text values in the specific order that reproduces
# the original author's shared-string table ordering (row 7's text was
# entered before row 6's text).
$ws.Cells.Item(2, 2).Value = "1,3,6"
$ws.Cells.Item(3, 2).Value = "1,3,7"
$ws.Cells.Item(4, 2).Value = "1,2,3,4,6,5,7"
$ws.Cells.Item(5, 2).Value = "1,2,9,12"
$ws.Cells.Item(7, 2).Value = "1,11,13,8,7"
$ws.Cells.Item(6, 2).Value = "1,10,13,8,7"
$ws.Cells.Item(8, 2).Value = "1,11,8,7"

# id column (A) values
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(8, 1).Value = 7

# Center-align the id/path data block (A2:B8) and the extended, otherwise
# empty, formatted region C6:D18 (columns C/D only gain formatting from row 6
# down, matching the original edit).
$ws.Range("A2:B8").HorizontalAlignment = -4108  # xlCenter
$ws.Range("C6:D18").HorizontalAlignment = -4108  # xlCenter

# Column B width (~34.27 characters)
$ws.Columns.Item(2).ColumnWidth = 33.5

# Selection
$ws.Range("K18").Select()
